$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 4 -> Correspond Handoff Datetime (D4) and Correspond Handback DateTime (G4)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D4").Value = "2016-01-25 05:48:32"
$wsZh.Range("G4").Value = "2016-01-25 05:49:19"

# de-de sheet: row 4 -> Correspond Handoff Datetime (D4) and Correspond Handback DateTime (G4)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D4").Value = "2016-01-25 05:48:43"
$wsDe.Range("G4").Value = "2016-01-25 05:49:36"
